$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10767.821
$ws.Range("I21").Value = 5894.684
$ws.Range("J21").Value = 21055.555
$ws.Range("K21").Value = 5894.684
$ws.Range("L21").Value = 21055.555
$ws.Range("M21").Value = -5426.684
$ws.Range("N21").Value = -21991.555
$ws.Range("H23").Value = 10767.821
$ws.Range("I23").Value = 5894.684
$ws.Range("J23").Value = 21055.555
$ws.Range("K23").Value = 5894.684
$ws.Range("L23").Value = 21055.555
$ws.Range("M23").Value = -5660.684
$ws.Range("N23").Value = -21523.555
$ws.Range("H93").Value = 25480.77
$ws.Range("J93").Value = 25480.77
$ws.Range("L93").Value = 25480.77
$ws.Range("N93").Value = -30472.77
$ws.Range("H94").Value = 1632.6666
$ws.Range("I94").Value = 1632.6666
$ws.Range("K94").Value = 1632.6666
$ws.Range("M94").Value = -1181.6666
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H109").Value = 33205.883
$ws.Range("J109").Value = 33205.883
$ws.Range("L109").Value = 33205.883
$ws.Range("N109").Value = -35979.883
$ws.Range("H116").Value = 308309.06
$ws.Range("J116").Value = 7768
$ws.Range("L116").Value = 7768
$ws.Range("N116").Value = -14652
$ws.Range("H129").Value = 843.8
$ws.Range("J129").Value = 874.3511
$ws.Range("L129").Value = 2623.0533
$ws.Range("N129").Value = -12623.0533
$ws.Range("H137").Value = 3665993
$ws.Range("I137").Value = 4763791
$ws.Range("J137").Value = 6666.6665
$ws.Range("K137").Value = 14291373
$ws.Range("L137").Value = 19999.9995
$ws.Range("M137").Value = -14288823
$ws.Range("N137").Value = -25099.9995
$ws.Range("H138").Value = 2559.27
$ws.Range("I138").Value = 700.8946999999999
$ws.Range("J138").Value = 2995.1853
$ws.Range("K138").Value = 2102.6841
$ws.Range("L138").Value = 8985.555899999999
$ws.Range("M138").Value = 3037.3159
$ws.Range("N138").Value = -19265.5559
$ws.Range("H141").Value = 24466.334
$ws.Range("I141").Value = 26899.625
$ws.Range("K141").Value = 80698.875
$ws.Range("M141").Value = -75518.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 911
$ws.Range("I2").Value = 911
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 911
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -798
$ws.Range("N2").ClearContents()
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H22").Value = 2666.6667
$ws.Range("I22").Value = 2666.6667
$ws.Range("K22").Value = 2666.6667
$ws.Range("M22").Value = -2367.6667
$ws.Range("H32").Value = 8017.9697
$ws.Range("I32").Value = 6689.3706
$ws.Range("J32").Value = 13996.667
$ws.Range("K32").Value = 6689.3706
$ws.Range("L32").Value = 13996.667
$ws.Range("M32").Value = -6402.3706
$ws.Range("N32").Value = -14570.667
$ws.Range("H45").Value = 2550
$ws.Range("I45").Value = 2400
$ws.Range("K45").Value = 2400
$ws.Range("M45").Value = -2023
$ws.Range("H74").Value = 9941
$ws.Range("I74").Value = 17031.8
$ws.Range("K74").Value = 17031.8
$ws.Range("M74").Value = -16157.8
$ws.Range("H77").Value = 9941
$ws.Range("I77").Value = 17031.8
$ws.Range("K77").Value = 85159
$ws.Range("M77").Value = -80791
$ws.Range("H92").Value = 30400
$ws.Range("J92").Value = 30400
$ws.Range("L92").Value = 30400
$ws.Range("N92").Value = -35392
$ws.Range("H116").Value = 911
$ws.Range("I116").Value = 911
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 911
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1383
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 1928.9231
$ws.Range("I132").Value = 1160.8636
$ws.Range("J132").Value = 6153.25
$ws.Range("K132").Value = 3482.5908
$ws.Range("L132").Value = 18459.75
$ws.Range("M132").Value = -952.5907999999999
$ws.Range("N132").Value = -23519.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 911
$ws.Range("I3").Value = 911
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 911
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -797
$ws.Range("N3").ClearContents()
$ws.Range("H92").Value = 62800
$ws.Range("J92").Value = 62800
$ws.Range("L92").Value = 62800
$ws.Range("N92").Value = -67792
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 346.8
$ws.Range("I7").Value = 359.66666
$ws.Range("J7").Value = 336.27274
$ws.Range("K7").Value = 359.66666
$ws.Range("L7").Value = 336.27274
$ws.Range("M7").Value = -246.66666
$ws.Range("N7").Value = -562.27274
$ws.Range("H31").Value = 4253.5713
$ws.Range("I31").Value = 1281.7778
$ws.Range("J31").Value = 9602.799999999999
$ws.Range("K31").Value = 1281.7778
$ws.Range("L31").Value = 9602.799999999999
$ws.Range("M31").Value = -986.7778000000001
$ws.Range("N31").Value = -10192.8
$ws.Range("H34").Value = 4253.5713
$ws.Range("I34").Value = 1281.7778
$ws.Range("J34").Value = 9602.799999999999
$ws.Range("K34").Value = 1281.7778
$ws.Range("L34").Value = 9602.799999999999
$ws.Range("M34").Value = -1079.7778
$ws.Range("N34").Value = -10006.8
$ws.Range("H107").Value = 861.3333
$ws.Range("I107").Value = 832.2
$ws.Range("J107").Value = 882.1429000000001
$ws.Range("K107").Value = 832.2
$ws.Range("L107").Value = 882.1429000000001
$ws.Range("M107").Value = 1087.8
$ws.Range("N107").Value = -4722.1429
$ws.Range("H137").Value = 42028.75
$ws.Range("J137").Value = 42028.75
$ws.Range("L137").Value = 42028.75
$ws.Range("N137").Value = -52228.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 20726.834
$ws.Range("I11").Value = 20726.834
$ws.Range("K11").Value = 62180.50199999999
$ws.Range("M11").Value = -62040.50199999999
$ws.Range("H75").Value = 20000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 20000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 60000
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -61996
$ws.Range("H78").Value = 20000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 20000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 180000
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -189984
$ws.Range("H107").Value = 486.88235
$ws.Range("I107").Value = 390
$ws.Range("J107").Value = 595.875
$ws.Range("K107").Value = 1170
$ws.Range("L107").Value = 1787.625
$ws.Range("M107").Value = 750
$ws.Range("N107").Value = -5627.625
$ws.Range("H113").Value = 3788475.8
$ws.Range("I113").Value = 624.13336
$ws.Range("J113").Value = 6945019
$ws.Range("K113").Value = 1872.40008
$ws.Range("L113").Value = 20835057
$ws.Range("M113").Value = 297.5999199999999
$ws.Range("N113").Value = -20839397
$ws.Range("H114").Value = 66668496
$ws.Range("I114").Value = 333333470
$ws.Range("J114").Value = 2255.4167
$ws.Range("K114").Value = 1000000410
$ws.Range("L114").Value = 6766.250100000001
$ws.Range("M114").Value = -999997156
$ws.Range("N114").Value = -13274.2501
$ws.Range("H126").Value = 2800
$ws.Range("I126").Value = 2800
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8400
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3460
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 667.65
$ws.Range("I131").Value = 248.32
$ws.Range("J131").Value = 807.4267
$ws.Range("K131").Value = 744.96
$ws.Range("L131").Value = 2422.2801
$ws.Range("M131").Value = 4295.04
$ws.Range("N131").Value = -12502.2801
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H115").Value = 39800
$ws.Range("J115").Value = 39800
$ws.Range("L115").Value = 39800
$ws.Range("N115").Value = -42150
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H69").Value = 18900
$ws.Range("J69").Value = 18900
$ws.Range("L69").Value = 18900
$ws.Range("N69").Value = -20398
$ws.Range("H72").Value = 18900
$ws.Range("J72").Value = 18900
$ws.Range("L72").Value = 56700
$ws.Range("N72").Value = -64188
$ws.Range("H132").Value = 9809187
$ws.Range("I132").Value = 5654.2085
$ws.Range("K132").Value = 16962.6255
$ws.Range("L132").Value = 100012998
$ws.Range("M132").Value = -14432.6255
$ws.Range("N132").Value = -100018058
